$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '60.247.85'
$ws.Cells.Item(2, 5).Value = '  -4.22%  '

$ws.Cells.Item(3, 4).Value = '2.986.00'
$ws.Cells.Item(3, 5).Value = '  -5.81%  '

$ws.Cells.Item(4, 5).Value = '  -0.03%  '

$ws.Cells.Item(5, 4).Value = '''578.73'
$ws.Cells.Item(5, 5).Value = '  -2.04%  '

$ws.Cells.Item(6, 4).Value = '''125.73'
$ws.Cells.Item(6, 5).Value = '  -6.54%  '

$ws.Cells.Item(7, 5).Value = '  +0.12%  '

$ws.Cells.Item(8, 4).Value = '2.981.61'
$ws.Cells.Item(8, 5).Value = '  -5.86%  '

$ws.Cells.Item(9, 5).Value = '  -2.61%  '

$ws.Cells.Item(10, 5).Value = '  -5.72%  '

$ws.Cells.Item(11, 5).Value = '  -2.13%  '

$ws.Cells.Item(12, 5).Value = '  -2.44%  '

$ws.Cells.Item(13, 4).Value = '''0.0000223'
$ws.Cells.Item(13, 5).Value = '  -5.57%  '

$ws.Cells.Item(14, 4).Value = '''32.54'
$ws.Cells.Item(14, 5).Value = '  -5.07%  '

$ws.Cells.Item(15, 5).Value = '  +0.19%  '

$ws.Cells.Item(16, 4).Value = '3.470.91'
$ws.Cells.Item(16, 5).Value = '  -5.95%  '

$ws.Cells.Item(17, 2).Value = 'WrappedEther'
$ws.Cells.Item(17, 3).Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Cells.Item(17, 4).Value = '2.984.66'
$ws.Cells.Item(17, 5).Value = '  -5.97%  '

$ws.Cells.Item(18, 2).Value = 'WrappedBTC'
$ws.Cells.Item(18, 3).Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Cells.Item(18, 4).Value = '60.157.98'
$ws.Cells.Item(18, 5).Value = '  -4.33%  '

$ws.Cells.Item(19, 4).Value = '''6.22'
$ws.Cells.Item(19, 5).Value = '  -4.71%  '

$ws.Cells.Item(20, 4).Value = '''432.28'
$ws.Cells.Item(20, 5).Value = '  -5.83%  '

$ws.Cells.Item(21, 4).Value = '''13.11'
$ws.Cells.Item(21, 5).Value = '  -6.08%  '

$ws.Cells.Item(22, 5).Value = '  -4.71%  '

$ws.Cells.Item(23, 4).Value = '''7.02'
$ws.Cells.Item(23, 5).Value = '  -7.29%  '

$ws.Cells.Item(24, 4).Value = '''12.67'
$ws.Cells.Item(24, 5).Value = '  -4.44%  '

$ws.Cells.Item(25, 4).Value = '''79.10'
$ws.Cells.Item(25, 5).Value = '  -3.84%  '

$ws.Cells.Item(26, 5).Value = '  +0.14%  '

$ws.Cells.Item(27, 4).Value = '''0.999'
$ws.Cells.Item(27, 5).Value = '  -0.05%  '

$ws.Cells.Item(28, 5).Value = '  -4.39%  '

$ws.Cells.Item(29, 4).Value = '''7.28'
$ws.Cells.Item(29, 5).Value = '  -4.37%  '

$ws.Cells.Item(30, 4).Value = '''1.89'
$ws.Cells.Item(30, 5).Value = '  -6.69%  '

$ws.Cells.Item(31, 5).Value = '  -8.82%  '

$ws.Cells.Item(32, 4).Value = '''25.32'
$ws.Cells.Item(32, 5).Value = '  -6.57%  '

$ws.Cells.Item(33, 4).Value = '''0.0935'
$ws.Cells.Item(33, 5).Value = '  -7.38%  '

$ws.Cells.Item(34, 4).Value = '''2.17'
$ws.Cells.Item(34, 5).Value = '  -7.89%  '

$ws.Cells.Item(35, 4).Value = '''0.951'
$ws.Cells.Item(35, 5).Value = '  -7.20%  '

$ws.Cells.Item(36, 4).Value = '''5.60'
$ws.Cells.Item(36, 5).Value = '  -3.09%  '

$ws.Cells.Item(37, 4).Value = '''49.52'
$ws.Cells.Item(37, 5).Value = '  -3.29%  '

$ws.Cells.Item(38, 5).Value = '  -6.38%  '

$ws.Cells.Item(39, 2).Value = 'Cosmos'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(39, 4).Value = '''7.99'
$ws.Cells.Item(39, 5).Value = '  -1.02%  '

$ws.Cells.Item(40, 2).Value = 'VeChain'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(40, 4).Value = '''0.0359'
$ws.Cells.Item(40, 5).Value = '  -6.80%  '

$ws.Cells.Item(41, 4).Value = '''386.13'
$ws.Cells.Item(41, 5).Value = '  -3.99%  '

$ws.Cells.Item(42, 5).Value = '  -2.36%  '

$ws.Cells.Item(43, 5).Value = '  -6.63%  '

$ws.Cells.Item(44, 4).Value = '2.629.24'

$ws.Cells.Item(45, 5).Value = '  +0.07%  '

$ws.Cells.Item(46, 4).Value = '''0.235'

$ws.Cells.Item(47, 4).Value = '''119.57'
$ws.Cells.Item(47, 5).Value = '  -3.72%  '

$ws.Cells.Item(48, 4).Value = '''1.99'
$ws.Cells.Item(48, 5).Value = '  -5.38%  '

$ws.Cells.Item(49, 5).Value = '  -3.52%  '

$ws.Cells.Item(50, 4).Value = '''23.45'
$ws.Cells.Item(50, 5).Value = '  -6.26%  '

$ws.Cells.Item(51, 4).Value = '''31.17'
$ws.Cells.Item(51, 5).Value = '  -9.88%  '
